# 13 Jul 2016:  All HTTP Calls have proper, functioning erorr handlers
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5: Status -> Closed, Short Desc -> updated text, clear the Long Desc
$ws.Range("C5").Value = "Update CB Functions / Error Handling / Front End"
$ws.Range("B5").Value = "Closed"
$ws.Range("D5").Clear()

# Add a new row 13 duplicating the new issue entry with status Open
$ws.Range("B13").Value = "Open"
$ws.Range("C13").Value = "Update CB Functions / Error Handling / Front End"
$ws.Range("B13:C13").EntireRow.RowHeight = 42

# Move the active selection to B14, matching the author's final cursor position
$ws.Range("B14").Select()
